$changes = @(
    @{ Row = 2; D = '25.969.97'; E = '  +0.51%  ' }
    @{ Row = 3; D = '1.650.39'; E = '  +1.03%  ' }
    @{ Row = 4; D = $null; E = '  +0.65%  ' }
    @{ Row = 5; D = '216.16'; E = '  +0.66%  ' }
    @{ Row = 6; D = '0.5122'; E = '  +2.17%  ' }
    @{ Row = 7; D = '1.007'; E = '  +0.52%  ' }
    @{ Row = 8; D = '0.2588'; E = '  +0.90%  ' }
    @{ Row = 9; D = '0.06439'; E = '  +0.57%  ' }
    @{ Row = 10; D = '19.76'; E = '  +0.91%  ' }
    @{ Row = 11; D = '0.07780'; E = '  +1.09%  ' }
    @{ Row = 12; D = '4.330'; E = '  +2.19%  ' }
    @{ Row = 13; D = '1.645.33'; E = '  +0.70%  ' }
    @{ Row = 14; D = '0.5491'; E = '  +1.30%  ' }
    @{ Row = 15; D = '0.0₅7912'; E = '  -0.06%  ' }
    @{ Row = 16; D = '65.06'; E = '  +2.67%  ' }
    @{ Row = 17; D = '26.045.98'; E = '  +0.78%  ' }
    @{ Row = 18; D = '1.007'; E = '  +0.55%  ' }
    @{ Row = 19; D = '199.22'; E = '  -1.32%  ' }
    @{ Row = 20; D = '4.471'; E = '  +3.45%  ' }
    @{ Row = 21; D = '10.07'; E = '  +1.49%  ' }
    @{ Row = 22; D = '6.079'; E = '  +1.90%  ' }
    @{ Row = 23; D = '1.009'; E = '  +0.62%  ' }
    @{ Row = 24; D = '1.868'; E = '  -2.14%  ' }
    @{ Row = 25; D = '140.51'; E = '  -0.19%  ' }
    @{ Row = 26; D = '0.1154'; E = '  +1.50%  ' }
    @{ Row = 27; D = '6.926'; E = '  +3.48%  ' }
    @{ Row = 28; D = '15.80'; E = '  +0.88%  ' }
    @{ Row = 29; D = '1.243'; E = '  +0.38%  ' }
    @{ Row = 30; D = '0.05038'; E = '  +1.36%  ' }
    @{ Row = 31; D = '3.298'; E = '  +1.42%  ' }
    @{ Row = 32; D = '3.216'; E = '  +1.38%  ' }
    @{ Row = 33; D = '1.549'; E = '  +0.82%  ' }
    @{ Row = 34; D = '2.370'; E = '  +0.36%  ' }
    @{ Row = 35; D = '0.8985'; E = '  +0.82%  ' }
    @{ Row = 36; D = '2.592'; E = '  -0.87%  ' }
    @{ Row = 37; D = '1.139.45'; E = '  -2.27%  ' }
    @{ Row = 38; D = '0.5566'; E = '  -0.52%  ' }
    @{ Row = 39; D = '0.01567'; E = $null }
    @{ Row = 40; D = '1.008'; E = '  +0.66%  ' }
    @{ Row = 41; D = '5.684'; E = '  +0.20%  ' }
    @{ Row = 42; D = '0.8190'; E = '  +1.56%  ' }
    @{ Row = 43; D = '100.03'; E = '  +0.83%  ' }
    @{ Row = 44; D = '0.0₈125'; E = '  +9.20%  ' }
    @{ Row = 45; D = '1.785.27'; E = '  +0.81%  ' }
    @{ Row = 46; D = '0.4540'; E = '  +0.56%  ' }
    @{ Row = 47; D = '55.51'; E = '  +1.55%  ' }
    @{ Row = 48; D = '1.008'; E = '  +0.62%  ' }
    @{ Row = 49; D = $null; E = '  +0.50%  ' }
    @{ Row = 50; D = '0.09604'; E = '  +3.75%  ' }
    @{ Row = 51; D = '1.006'; E = '  +0.47%  ' }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($change in $changes) {
    $row = $change.Row
    if ($null -ne $change.D) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $change.D
        $cell.Style = "Normal"
    }
    if ($null -ne $change.E) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $change.E
        $cell.Style = "Normal"
    }
}
